$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "LOCAL COMERCIAL" -> "NEGOCIOS" as a new row 6 (pushes existing rows down)
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "LOCAL COMERCIAL"
$ws.Range("B6").Value = "NEGOCIOS"

# Insert "OFICINAS" -> "LUGARES DE TRABAJO" as a new row 24 (pushes existing rows down)
$ws.Rows.Item(24).Insert()
$ws.Range("A24").Value = "OFICINAS"
$ws.Range("B24").Value = "LUGARES DE TRABAJO"

# Update the view: select B25 (also clears any stale scroll position)
$ws.Range("B25").Select() | Out-Null
